$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EV Home win")
$ws.Rows.Item(16).Delete()
